# Generate Report for Handback
# Refresh the handoff/handback timestamps recorded for the
# 96c4c545-2127-42e8-a97e-09db01a99ce8 file (row 2 in each report sheet)
# to reflect the newly regenerated XLIFF handback.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G2").Value = "2016-09-05 00:55:11"

# --- zh-cn sheet ------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-09-05 00:55:00"
$zhcn.Range("K2").Value = "2016-09-05 00:55:29"

# --- de-de sheet ------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H2").Value = "2016-09-05 00:55:11"
$dede.Range("K2").Value = "2016-09-05 00:55:37"
